# Insert 3 new data rows before the existing row 479 (shifting the old
# rows 479-567 down to 482-570), then populate the 3 new rows with the
# new data points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at position 479 (existing content shifts down).
$ws.Rows("479:481").Insert()

# Common / unchanged values for this homogeneous table.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$catId     = 100112033
$categoria = "Lechuga"
$calidad   = "Primera"
$clasif    = "Hortaliza"

# New row 479
$ws.Range("A479").Value = $mercadoId
$ws.Range("B479").Value = $mercado
$ws.Range("C479").Value = $region
$ws.Range("D479").Value = 44617
$ws.Range("E479").Value = $codreg
$ws.Range("F479").Value = $catId
$ws.Range("G479").Value = $categoria
$ws.Range("H479").Value = "Conconina(o)"
$ws.Range("I479").Value = $calidad
$ws.Range("J479").Value = 240
$ws.Range("K479").Value = 5500
$ws.Range("L479").Value = 6000
$ws.Range("M479").Value = 5750
$ws.Range("N479").Value = "$/caja 10 unidades"
$ws.Range("O479").Value = "Provincia de Diguillín"
$ws.Range("P479").Value = 575
$ws.Range("Q479").Value = 10
$ws.Range("R479").Value = $clasif

# New row 480
$ws.Range("A480").Value = $mercadoId
$ws.Range("B480").Value = $mercado
$ws.Range("C480").Value = $region
$ws.Range("D480").Value = 44617
$ws.Range("E480").Value = $codreg
$ws.Range("F480").Value = $catId
$ws.Range("G480").Value = $categoria
$ws.Range("H480").Value = "Escarola"
$ws.Range("I480").Value = $calidad
$ws.Range("J480").Value = 200
$ws.Range("K480").Value = 7000
$ws.Range("L480").Value = 7500
$ws.Range("M480").Value = 7250
$ws.Range("N480").Value = "$/caja 15 unidades"
$ws.Range("O480").Value = "Región del Maule"
$ws.Range("P480").Value = 483
$ws.Range("Q480").Value = 15
$ws.Range("R480").Value = $clasif

# New row 481
$ws.Range("A481").Value = $mercadoId
$ws.Range("B481").Value = $mercado
$ws.Range("C481").Value = $region
$ws.Range("D481").Value = 44617
$ws.Range("E481").Value = $codreg
$ws.Range("F481").Value = $catId
$ws.Range("G481").Value = $categoria
$ws.Range("H481").Value = "Marina"
$ws.Range("I481").Value = $calidad
$ws.Range("J481").Value = 200
$ws.Range("K481").Value = 5500
$ws.Range("L481").Value = 6000
$ws.Range("M481").Value = 5750
$ws.Range("N481").Value = "$/caja 18 unidades"
$ws.Range("O481").Value = "Región del Maule"
$ws.Range("P481").Value = 319
$ws.Range("Q481").Value = 18
$ws.Range("R481").Value = $clasif

# Ensure the date column keeps its date/time number format for the
# newly inserted rows (style index 2 in the original workbook).
$ws.Range("D479:D481").NumberFormat = "YYYY-MM-DD HH:MM:SS"
